# Update "想去人数" (F column) values on the "展览", "演出" and "全部类型" sheets
# to match the newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" -----------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$updates1 = @{
    2  = 3163
    3  = 566
    4  = 1138
    5  = 123
    6  = 87
    8  = 61
    9  = 1171
    10 = 16504
    11 = 284
    13 = 1043
    14 = 6403
    15 = 643
    16 = 133
    17 = 83
    19 = 128
    21 = 53
    26 = 6
    29 = 902
    30 = 62
    33 = 11401
    35 = 21
    36 = 156
    38 = 3848
}
foreach ($row in $updates1.Keys) {
    $ws1.Range("F$row").Value = $updates1[$row]
}

# --- Sheet "演出" -------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 22

# --- Sheet "全部类型" ----------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$updates4 = @{
    2  = 3163
    3  = 566
    4  = 1138
    5  = 123
    6  = 87
    8  = 61
    9  = 1171
    10 = 16504
    11 = 284
    13 = 1043
    14 = 6403
    15 = 643
    16 = 133
    17 = 83
    19 = 128
    21 = 53
    26 = 6
    29 = 902
    30 = 62
    33 = 22
    34 = 11401
    36 = 21
    37 = 156
    39 = 3848
}
foreach ($row in $updates4.Keys) {
    $ws4.Range("F$row").Value = $updates4[$row]
}

$wb.Save()
